{"js": "\n// Mapping: the table has 20 rows x 5 columns = 100 cells, in row-major\n// order matching the order the math-fact cells appear in the document.\n// Each cell's existing text (e.g. \"85-17=68\") is replaced with the new\n// value (e.g. \"96-62=34\") per the commit's diff, preserving all\n// paragraph/run formatting (font, size) since we only change the text.\nconst NEW_VALUES = [\"96-62=34\", \"48-8=40\", \"40+16=56\", \"32-12=20\", \"45+36=81\", \"32-26=6\", \"38+37=75\", \"31+31=62\", \"37+19=56\", \"6+5=11\", \"90+2=92\", \"27+37=64\", \"15+23=38\", \"64-17=47\", \"69+14=83\", \"56-55=1\", \"6+54=60\", \"34+13=47\", \"45+12=57\", \"6+78=84\", \"41+57=98\", \"9+83=92\", \"43+43=86\", \"0+91=91\", \"81-46=35\", \"60-21=39\", \"65+20=85\", \"14+45=59\", \"90-36=54\", \"93-2=91\", \"23+40=63\", \"89-77=12\", \"33+61=94\", \"44-35=9\", \"15+66=81\", \"5+57=62\", \"59-26=33\", \"39-30=9\", \"27+7=34\", \"74-63=11\", \"92-8=84\", \"94-65=29\", \"7-4=3\", \"85-14=71\", \"20+71=91\", \"13-0=13\", \"10+0=10\", \"57-27=30\", \"47+19=66\", \"19+11=30\", \"50-4=46\", \"57-53=4\", \"99-38=61\", \"62-6=56\", \"78-1=77\", \"41-23=18\", \"39+34=73\", \"38-28=10\", \"24+6=30\", \"92-71=21\", \"75-65=10\", \"15-6=9\", \"25+50=75\", \"89-15=74\", \"1+40=41\", \"32+5=37\", \"74-44=30\", \"14-6=8\", \"81+0=81\", \"69-2=67\", \"54-31=23\", \"32+51=83\", \"83-33=50\", \"69+17=86\", \"1+96=97\", \"7+54=61\", \"63-32=31\", \"74-4=70\", \"74+10=84\", \"1+2=3\", \"6+43=49\", \"98-45=53\", \"49-0=49\", \"23+27=50\", \"48-22=26\", \"36+13=49\", \"55-7=48\", \"74-65=9\", \"81-58=23\", \"15-12=3\", \"16+2=18\", \"48+1=49\", \"56-48=8\", \"41-20=21\", \"55-9=46\", \"38+55=93\", \"65-58=7\", \"17+67=84\", \"79-57=22\", \"10+26=36\"];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document, found none.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load cells for every row up front.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    const para = cell.body.paragraphs.getFirst();\n    para.load(\"text\");\n    // Defer sync until batched below would be nicer, but keep it simple &\n    // robust: sync per cell paragraph load then write.\n    await context.sync();\n\n    if (idx < NEW_VALUES.length) {\n      // Replace the run text in-place (keeps run/paragraph formatting such\n      // as font TimeNewRoman / size 30 and the cell's paragraph alignment).\n      const range = para.getRange();\n      range.insertText(NEW_VALUES[idx], Word.InsertLocation.replace);\n    }\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each math-fact cell's text with its updated value, per the\n# commit diff. The document contains a single 20x5 table of cells like\n# \"85-17=68\"; every one of the 100 cells is rewritten to a new expression\n# (e.g. \"96-62=34\"), in the same order they appear in the document. Since\n# none of the old values repeat and none of the new values collide with a\n# not-yet-processed old value, a simple ordered Find/Replace (one hit at a\n# time, starting the search from the top of the document each time) is\n# unambiguous and safe.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"85-17=68\"; New = \"96-62=34\" },\n    @{ Old = \"15+30=45\"; New = \"48-8=40\" },\n    @{ Old = \"11+70=81\"; New = \"40+16=56\" },\n    @{ Old = \"20+67=87\"; New = \"32-12=20\" },\n    @{ Old = \"73-30=43\"; New = \"45+36=81\" },\n    @{ Old = \"65-36=29\"; New = \"32-26=6\" },\n    @{ Old = \"33+20=53\"; New = \"38+37=75\" },\n    @{ Old = \"88-79=9\"; New = \"31+31=62\" },\n    @{ Old = \"8+90=98\"; New = \"37+19=56\" },\n    @{ Old = \"9-4=5\"; New = \"6+5=11\" },\n    @{ Old = \"84+8=92\"; New = \"90+2=92\" },\n    @{ Old = \"56+19=75\"; New = \"27+37=64\" },\n    @{ Old = \"18+58=76\"; New = \"15+23=38\" },\n    @{ Old = \"26+38=64\"; New = \"64-17=47\" },\n    @{ Old = \"70-63=7\"; New = \"69+14=83\" },\n    @{ Old = \"72+22=94\"; New = \"56-55=1\" },\n    @{ Old = \"78-36=42\"; New = \"6+54=60\" },\n    @{ Old = \"58-34=24\"; New = \"34+13=47\" },\n    @{ Old = \"78-12=66\"; New = \"45+12=57\" },\n    @{ Old = \"38-3=35\"; New = \"6+78=84\" },\n    @{ Old = \"70-4=66\"; New = \"41+57=98\" },\n    @{ Old = \"51-25=26\"; New = \"9+83=92\" },\n    @{ Old = \"19+12=31\"; New = \"43+43=86\" },\n    @{ Old = \"43-2=41\"; New = \"0+91=91\" },\n    @{ Old = \"61-27=34\"; New = \"81-46=35\" },\n    @{ Old = \"59+4=63\"; New = \"60-21=39\" },\n    @{ Old = \"94-92=2\"; New = \"65+20=85\" },\n    @{ Old = \"14+16=30\"; New = \"14+45=59\" },\n    @{ Old = \"72-15=57\"; New = \"90-36=54\" },\n    @{ Old = \"34-28=6\"; New = \"93-2=91\" },\n    @{ Old = \"12+81=93\"; New = \"23+40=63\" },\n    @{ Old = \"45-32=13\"; New = \"89-77=12\" },\n    @{ Old = \"88+10=98\"; New = \"33+61=94\" },\n    @{ Old = \"78-21=57\"; New = \"44-35=9\" },\n    @{ Old = \"20+74=94\"; New = \"15+66=81\" },\n    @{ Old = \"18-5=13\"; New = \"5+57=62\" },\n    @{ Old = \"32+42=74\"; New = \"59-26=33\" },\n    @{ Old = \"47+0=47\"; New = \"39-30=9\" },\n    @{ Old = \"98-77=21\"; New = \"27+7=34\" },\n    @{ Old = \"7+56=63\"; New = \"74-63=11\" },\n    @{ Old = \"68-13=55\"; New = \"92-8=84\" },\n    @{ Old = \"60-0=60\"; New = \"94-65=29\" },\n    @{ Old = \"46+33=79\"; New = \"7-4=3\" },\n    @{ Old = \"23+71=94\"; New = \"85-14=71\" },\n    @{ Old = \"44-24=20\"; New = \"20+71=91\" },\n    @{ Old = \"73-3=70\"; New = \"13-0=13\" },\n    @{ Old = \"97-12=85\"; New = \"10+0=10\" },\n    @{ Old = \"7+9=16\"; New = \"57-27=30\" },\n    @{ Old = \"62-4=58\"; New = \"47+19=66\" },\n    @{ Old = \"53+41=94\"; New = \"19+11=30\" },\n    @{ Old = \"19+36=55\"; New = \"50-4=46\" },\n    @{ Old = \"83-38=45\"; New = \"57-53=4\" },\n    @{ Old = \"72-71=1\"; New = \"99-38=61\" },\n    @{ Old = \"77-29=48\"; New = \"62-6=56\" },\n    @{ Old = \"47+10=57\"; New = \"78-1=77\" },\n    @{ Old = \"86-4=82\"; New = \"41-23=18\" },\n    @{ Old = \"63-38=25\"; New = \"39+34=73\" },\n    @{ Old = \"98+0=98\"; New = \"38-28=10\" },\n    @{ Old = \"73+21=94\"; New = \"24+6=30\" },\n    @{ Old = \"96-6=90\"; New = \"92-71=21\" },\n    @{ Old = \"41+13=54\"; New = \"75-65=10\" },\n    @{ Old = \"99-81=18\"; New = \"15-6=9\" },\n    @{ Old = \"75-57=18\"; New = \"25+50=75\" },\n    @{ Old = \"86-42=44\"; New = \"89-15=74\" },\n    @{ Old = \"9-0=9\"; New = \"1+40=41\" },\n    @{ Old = \"49+4=53\"; New = \"32+5=37\" },\n    @{ Old = \"15+60=75\"; New = \"74-44=30\" },\n    @{ Old = \"99-95=4\"; New = \"14-6=8\" },\n    @{ Old = \"12+38=50\"; New = \"81+0=81\" },\n    @{ Old = \"27+53=80\"; New = \"69-2=67\" },\n    @{ Old = \"17+60=77\"; New = \"54-31=23\" },\n    @{ Old = \"4-4=0\"; New = \"32+51=83\" },\n    @{ Old = \"1+24=25\"; New = \"83-33=50\" },\n    @{ Old = \"26+50=76\"; New = \"69+17=86\" },\n    @{ Old = \"79-15=64\"; New = \"1+96=97\" },\n    @{ Old = \"10+51=61\"; New = \"7+54=61\" },\n    @{ Old = \"87-54=33\"; New = \"63-32=31\" },\n    @{ Old = \"38-10=28\"; New = \"74-4=70\" },\n    @{ Old = \"42-15=27\"; New = \"74+10=84\" },\n    @{ Old = \"60-22=38\"; New = \"1+2=3\" },\n    @{ Old = \"60-34=26\"; New = \"6+43=49\" },\n    @{ Old = \"33+9=42\"; New = \"98-45=53\" },\n    @{ Old = \"55-13=42\"; New = \"49-0=49\" },\n    @{ Old = \"33+40=73\"; New = \"23+27=50\" },\n    @{ Old = \"48+3=51\"; New = \"48-22=26\" },\n    @{ Old = \"26-20=6\"; New = \"36+13=49\" },\n    @{ Old = \"67-18=49\"; New = \"55-7=48\" },\n    @{ Old = \"44-25=19\"; New = \"74-65=9\" },\n    @{ Old = \"61+7=68\"; New = \"81-58=23\" },\n    @{ Old = \"58+17=75\"; New = \"15-12=3\" },\n    @{ Old = \"57+29=86\"; New = \"16+2=18\" },\n    @{ Old = \"86-69=17\"; New = \"48+1=49\" },\n    @{ Old = \"33+36=69\"; New = \"56-48=8\" },\n    @{ Old = \"15+81=96\"; New = \"41-20=21\" },\n    @{ Old = \"45+5=50\"; New = \"55-9=46\" },\n    @{ Old = \"64+10=74\"; New = \"38+55=93\" },\n    @{ Old = \"95-10=85\"; New = \"65-58=7\" },\n    @{ Old = \"26+31=57\"; New = \"17+67=84\" },\n    @{ Old = \"67-42=25\"; New = \"79-57=22\" },\n    @{ Old = \"86-0=86\"; New = \"10+26=36\" }\n)\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pair.New, $wdReplaceOne) | Out-Null\n}\n\nWrite-Output \"Replaced $($pairs.Count) cells\"\n"}
